{"js": "// Remove every paragraph in the document body except the first one\n// (\"Template Header\"). This deletes the \"Test Paper\" heading, the\n// intro paragraph, the figure (image + caption), the \"Conclusion\"\n// heading, and the closing paragraph \u2014 collapsing the document down\n// to just the template header, per the commit's cleanup.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (let i = paragraphs.items.length - 1; i >= 1; i--) {\n  paragraphs.items[i].delete();\n}\nawait context.sync();\n", "ps1": "# Remove every paragraph in the document except the first one\n# (\"Template Header\"). This deletes the \"Test Paper\" heading, the\n# intro paragraph, the figure (image + caption), the \"Conclusion\"\n# heading, and the closing paragraph \u2014 collapsing the document down\n# to just the template header, per the commit's cleanup.\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\nfor ($i = $count; $i -ge 2; $i--) {\n    $d.Paragraphs.Item($i).Range.Delete()\n}\n"}
